$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 158 (shifts existing rows 158-187 down to 159-188)
$ws.Rows.Item(158).Insert()

$row = 158
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44476
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112001
$ws.Cells.Item($row, 7).Value = "Berenjena"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 80
$ws.Cells.Item($row, 11).Value = 10000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 10000
$ws.Cells.Item($row, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 167
$ws.Cells.Item($row, 17).Value = 60
$ws.Cells.Item($row, 18).Value = "Hortaliza"
